# Adds two new "BRUNO DE FRAGA" test rows (19 and 20) to the Clientes sheet,
# mirroring the existing data pattern (column D - Endereco - stays empty).
# Values that look numeric must be forced to stay as text, like the rest of
# the sheet (Excel would otherwise silently convert them to numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 19
$ws.Range("A19").Value = "BRUNO DE FRAGA"
Set-TextValue "B19" "123123"
Set-TextValue "C19" "123123"
$ws.Range("E19").Value = "92320-195"
$ws.Range("F19").Value = "qweqweqweq@qweqwe"
Set-TextValue "G19" "123123123123"
$ws.Range("H19").Value = "Rua 3 Pinheiros I, 27"

# Row 20
$ws.Range("A20").Value = "BRUNO DE FRAGA"
Set-TextValue "B20" "123123"
Set-TextValue "C20" "123123"
$ws.Range("E20").Value = "92320-195"
$ws.Range("F20").Value = "1231!@3123"
Set-TextValue "G20" "123123"
$ws.Range("H20").Value = "Rua 3 Pinheiros I, 27"
